# Apply the "Add UML sequence diagrams content" edit:
#   - The whole admin/professor/student UML diagram (rectangles, connectors,
#     multiplicity textboxes) is shifted up and to the right on the slide
#     by a uniform amount (dx=+485800 EMU, dy=-2303512 EMU).
#   - The now-unused empty "Title" placeholder shape is removed.
#   - The presentation's custom first-slide-number override is cleared
#     (back to the default of 1).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Target absolute Left/Top (in points; 1 pt = 12700 EMU) for every shape
# that belongs to the UML diagram, taken from the new ppt/slides/slide1.xml.
$moves = @(
  @{ Name = "Rectangle 24";         Left = 266.251968503937;  Top = 203.53984251968504 },
  @{ Name = "Rectangle 25";         Left = 542.2519685039371; Top = 203.08047244094487 },
  @{ Name = "Elbow Connector 28";   Left = 368.251968503937;  Top = 217.62110236220474 },
  @{ Name = "TextBox 29";           Left = 476.251968503937;  Top = 194.64763779527559 },
  @{ Name = "TextBox 30";           Left = 368.251968503937;  Top = 191.53984251968504 },
  @{ Name = "Rectangle 31";         Left = 74.25196850393701; Top = 203.53984251968504 },
  @{ Name = "Elbow Connector 33";   Left = 356.02228346456695;Top = 1.3914173228346456 },
  @{ Name = "TextBox 34";           Left = 218.251968503937;  Top = 190.62110236220474 },
  @{ Name = "Elbow Connector 38";   Left = 176.251968503937;  Top = 216.01330708661416 },
  @{ Name = "TextBox 41";           Left = 512.2519685039371; Top = 232.62110236220474 },
  @{ Name = "TextBox 15";           Left = 176.251968503937;  Top = 195.7023622047244 },
  @{ Name = "TextBox 16";           Left = 104.25196850393701;Top = 245.7696062992126 },
  @{ Name = "TextBox 17";           Left = 590.2519685039371; Top = 232.62110236220474 }
)

foreach ($mv in $moves) {
  $shp = $s.Shapes.Item($mv.Name)
  $shp.Left = $mv.Left
  $shp.Top = $mv.Top
}

# Drop the empty title placeholder that used to sit above the diagram.
$title = $s.Shapes.Item("Title 35")
$title.Delete()

# The presentation no longer pins a custom starting slide number.
$p.PageSetup.FirstSlideNumber = 1
